{"js": "// Replace the date and each of the division problems in the table with\n// their updated values. Every source string below occurs exactly once in\n// the document, so a targeted search + replace keeps all run/paragraph\n// formatting intact.\nconst replacements = [\n  [\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"],\n  [\"62\u00f72=\", \"97\u00f74=\"],\n  [\"89\u00f77=\", \"96\u00f78=\"],\n  [\"53\u00f74=\", \"95\u00f79=\"],\n  [\"67\u00f73=\", \"54\u00f76=\"],\n  [\"68\u00f79=\", \"21\u00f78=\"],\n  [\"31\u00f74=\", \"27\u00f78=\"],\n  [\"64\u00f76=\", \"35\u00f75=\"],\n  [\"74\u00f74=\", \"14\u00f73=\"],\n  [\"53\u00f78=\", \"23\u00f77=\"],\n  [\"59\u00f79=\", \"35\u00f77=\"],\n  [\"22\u00f76=\", \"14\u00f78=\"],\n  [\"57\u00f77=\", \"62\u00f78=\"],\n  [\"46\u00f75=\", \"13\u00f77=\"],\n  [\"78\u00f72=\", \"19\u00f77=\"],\n  [\"20\u00f79=\", \"13\u00f74=\"],\n  [\"82\u00f79=\", \"65\u00f77=\"],\n  [\"58\u00f72=\", \"86\u00f76=\"],\n  [\"28\u00f78=\", \"96\u00f76=\"],\n  [\"29\u00f72=\", \"14\u00f75=\"],\n  [\"86\u00f73=\", \"32\u00f75=\"],\n  [\"77\u00f74=\", \"15\u00f76=\"],\n  [\"13\u00f79=\", \"37\u00f74=\"],\n  [\"71\u00f76=\", \"79\u00f76=\"],\n  [\"68\u00f75=\", \"67\u00f75=\"],\n  [\"42\u00f76=\", \"81\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each of the division problems in the table with\n# their updated values. Every source string below occurs exactly once in\n# the document, so Find/Replace on $d.Content keeps all run/paragraph\n# formatting intact (Find.Execute replaces text in place, run-by-run).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-17 Monday\", \"2024-06-18 Tuesday\"),\n    @(\"62\u00f72=\", \"97\u00f74=\"),\n    @(\"89\u00f77=\", \"96\u00f78=\"),\n    @(\"53\u00f74=\", \"95\u00f79=\"),\n    @(\"67\u00f73=\", \"54\u00f76=\"),\n    @(\"68\u00f79=\", \"21\u00f78=\"),\n    @(\"31\u00f74=\", \"27\u00f78=\"),\n    @(\"64\u00f76=\", \"35\u00f75=\"),\n    @(\"74\u00f74=\", \"14\u00f73=\"),\n    @(\"53\u00f78=\", \"23\u00f77=\"),\n    @(\"59\u00f79=\", \"35\u00f77=\"),\n    @(\"22\u00f76=\", \"14\u00f78=\"),\n    @(\"57\u00f77=\", \"62\u00f78=\"),\n    @(\"46\u00f75=\", \"13\u00f77=\"),\n    @(\"78\u00f72=\", \"19\u00f77=\"),\n    @(\"20\u00f79=\", \"13\u00f74=\"),\n    @(\"82\u00f79=\", \"65\u00f77=\"),\n    @(\"58\u00f72=\", \"86\u00f76=\"),\n    @(\"28\u00f78=\", \"96\u00f76=\"),\n    @(\"29\u00f72=\", \"14\u00f75=\"),\n    @(\"86\u00f73=\", \"32\u00f75=\"),\n    @(\"77\u00f74=\", \"15\u00f76=\"),\n    @(\"13\u00f79=\", \"37\u00f74=\"),\n    @(\"71\u00f76=\", \"79\u00f76=\"),\n    @(\"68\u00f75=\", \"67\u00f75=\"),\n    @(\"42\u00f76=\", \"81\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
